$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the rank/prolificid/name values between row 4 and row 5
# (re-ranked entries: "Annes" now ranks above "Sabrina")
$bTemp = $ws.Range("B4").Value2
$cTemp = $ws.Range("C4").Value2
$dTemp = $ws.Range("D4").Value2

$ws.Range("B4").Value2 = $ws.Range("B5").Value2
$ws.Range("C4").Value2 = $ws.Range("C5").Value2
$ws.Range("D4").Value2 = $ws.Range("D5").Value2

$ws.Range("B5").Value2 = $bTemp
$ws.Range("C5").Value2 = $cTemp
$ws.Range("D5").Value2 = $dTemp

# Updated matrices scores (recomputed)
$ws.Range("F2").Value2 = 14.09110926329862
$ws.Range("F3").Value2 = 13.45595228309568
$ws.Range("F4").Value2 = 13.21138131302576
$ws.Range("F5").Value2 = 13.05510090847672
$ws.Range("F6").Value2 = 12.11138011880338
$ws.Range("F7").Value2 = 10.38169935728711
$ws.Range("F8").Value2 = 10.3148949458874
$ws.Range("F9").Value2 = 8.206853693142603
$ws.Range("F10").Value2 = 5.433954494785023
$ws.Range("F11").Value2 = 2.475607795998219
$ws.Range("F12").Value2 = 2.127299389597505
$ws.Range("F13").Value2 = 1.091413571818724
